# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Each entry maps a 1-based row number to its new value for that sheet.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    2  = 6797
    3  = 15
    4  = 431
    7  = 546
    8  = 103
    9  = 98
    12 = 12
    13 = 175
    14 = 416
    16 = 1619
    17 = 23
    18 = 3417
    20 = 230
    21 = 9
    22 = 2048
    23 = 153
    24 = 2
    28 = 7
}

$sheet4Updates = @{
    2  = 6797
    3  = 15
    4  = 431
    8  = 546
    9  = 103
    10 = 98
    13 = 12
    14 = 175
    15 = 416
    17 = 1619
    18 = 23
    19 = 3417
    21 = 230
    22 = 9
    23 = 2048
    24 = 153
    25 = 2
    29 = 7
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
